$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 18 (week 12): hours logged increased from 2 to 4.5
$ws.Range("E18").Value = 4.5

# Row 18 (week 12): activity log text updated with bug-fix / likes-feature notes
$ws.Range("F18").Value = 'Added new "Likes" table wrote backend code to update likes in the database. Added backend code to close a session.. Fixed bugs in deleting message and fixed bugs in Junit tests involving the test to update message content and the test to delete messages. Fixed a JUnit test bug where we were testing the wrong procedure. Wrote SPROC to toggle likes on a message'

# Scroll the sheet view so row 11 is the top-left visible row (was row 9)
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 2
